$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "기타" (etc) response text is replaced everywhere by the fuller
# chatbot reply string. Once every cell referencing the old "기타" string
# is repointed, that shared string becomes unused and drops out of the
# workbook on save.
$ws.Range("C13").Value = "(기타)을(를) 맡기러 오셨군요?"
$ws.Range("C14").Value = "(기타)을(를) 맡기러 오셨군요?"
$ws.Range("C15").Value = "(기타)을(를) 맡기러 오셨군요?"
$ws.Range("C16").Value = "(기타)을(를) 맡기러 오셨군요?"

# D2 keeps the same displayed text ("classification") - just re-assert it.
$ws.Range("D2").Value = "classification"

# New training row for the "time" intent / slot-filling example.
$ws.Range("A17").Value = "11월 3일 10시 47분"
$ws.Range("B17").Value = "월|일|시|분"
$ws.Range("D17").Value = "time"

# Move the active selection, matching the author's last cursor position.
[void]$ws.Range("F4").Select()
